$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the date values for B3:D3 (Jan-2017, Feb-2017, Mar-2017) using
# the underlying serial date numbers so Excel doesn't auto-assign a
# built-in date number format before we apply the custom one.
$ws.Range("B3").Value = 42736
$ws.Range("C3").Value = 42767
$ws.Range("D3").Value = 42795

# Apply the custom date format
$ws.Range("B3:D3").NumberFormat = "mmm\-yyyy"

# Autofit columns B:D to match bestFit widths from the diff
$ws.Columns("B:D").AutoFit()

# Set the active selection to F2
$ws.Range("F2").Select()
